# Generate Report for Handback
#
# The localization-status report is refreshed after a handback run: the
# "b5467805-8898-4c89-9d0a-ffb3e3fe9350" file's status moves from
# "Ready for handoff" to "Handback transform failed", and the zh-cn / de-de
# detail sheets get an explanatory error message recorded in their
# "Error Detail" column (L) for that same file's row.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# Status column (shared across Overview/zh-cn/de-de) for the
# b5467805-8898-4c89-9d0a-ffb3e3fe9350 file row now reports a failed
# handback transform instead of "Ready for handoff".
$ws1.Range("B3").Value = "Handback transform failed"
$ws1.Range("C3").Value = "Handback transform failed"
$ws2.Range("C3").Value = "Handback transform failed"
$ws3.Range("C3").Value = "Handback transform failed"

# Error Detail (column L) explaining why the handback failed for each locale.
$ws2.Range("L3").Value = "Handback file name: 31xxdx50.4kh is different with handoff file name: b5467805-8898-4c89-9d0a-ffb3e3fe9350.56721758d12da70f7a393ec69e32b96f3c6a71cf.zh-cn."
$ws3.Range("L3").Value = "Handback file name: 31xxdx50.4kh is different with handoff file name: b5467805-8898-4c89-9d0a-ffb3e3fe9350.56721758d12da70f7a393ec69e32b96f3c6a71cf.de-de."
